# Updated cryptos list (Price / Volume(1h) columns) with refreshed market data.
# D<row> = Price text, E<row> = Volume(1h) text ("  +x.xx%  " / "  -x.xx%  ").
# Price cells whose text parses as a plain number get NumberFormat "@" first
# so Excel stores them as text (matching the original inline-string values)
# instead of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.600.62"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.821.15"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.31"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4662"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3588"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07118"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8991"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07800"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "1.849.48"
$ws.Range("E13").Value = "  +4.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.248"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.14"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008532"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "26.644.92"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.004"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.53"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.935"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.97"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.968"
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.57"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08788"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.123"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.735"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7279"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.430"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.074"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01924"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.920"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05097"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.812"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5020"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1490"
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.959"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.010"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4634"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.962"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.37"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.550"
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05993"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.60"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.63"
$ws.Range("E51").Value = "  -1.63%  "
